# Jessie Jr. Alcuizar Hellera (Q0159) - Training Dashboard update
# - header row (row 2) gets a white, bold font (color fix on the existing bold style)
# - H3 "PERIOD TO EXPIRE" recalculated: -42 -> -50
# - I3 "LAST UPDATE" bumped: 08-Sep-2025 -> 16-Sep-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Header row: make the (already bold) header font white.
# The header style is shared by both sheets, so update it everywhere it's used.
$headerRange = $ws.Range("A2:K2")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215   # RGB(255,255,255) -> white

$ws2 = $wb.Worksheets.Item("Exam Dashboard")
$headerRange2 = $ws2.Range("A2:G2")
$headerRange2.Font.Bold = $true
$headerRange2.Font.Color = 16777215   # RGB(255,255,255) -> white

# H3: period-to-expire value
$ws.Range("H3").Value = -50

# I3: last-update date, kept as literal text (not an Excel date serial).
# A leading apostrophe forces text-entry semantics like typing '16-Sep-2025 in the UI.
$ws.Range("I3").Value = "'16-Sep-2025"
